# Add an "OBJECTIVE:" paragraph (styled like the other section headings)
# right before the "EDUCATION:" heading, and move the stray "_GoBack"
# bookmark (previously sitting at the very end of the document, after
# "Fluent in Spanish") into the middle of the new objective sentence -
# which is exactly where Word leaves it after a user types new text and
# saves.

$d = $word.ActiveDocument

# The document currently has a left-over "_GoBack" bookmark near the end
# of the doc (after "...Fluent in Spanish"). Remove it - it will be
# re-created at its new location below. (Bookmark names must be unique,
# so adding the new one further down would implicitly delete this one
# anyway; we do it explicitly for clarity.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Split the "EDUCATION:" paragraph so that a brand new paragraph (which
# inherits the same Heading1 style) is created immediately before it,
# containing our objective text.
$replacement = "OBJECTIVE: Eager to drive solutions at American Express on a full-time basis^pEDUCATION:"
$d.Content.Find.Execute("EDUCATION:", $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null

# Locate the freshly created "OBJECTIVE: ..." paragraph (now the
# paragraph immediately before "EDUCATION:").
$objPara = $d.Paragraphs.Item(4)
$pStart = $objPara.Range.Start

# Make sure the new paragraph uses the same heading style as the rest of
# the section headers.
$objPara.Style = "Heading 1"

# Break the sentence following "OBJECTIVE: " into the same run chunks
# used in the final document, and give each of them the smaller 12pt
# (sz=24) font size used for the objective blurb (the "OBJECTIVE: "
# label itself keeps the default Heading1 size).
$prefix = "OBJECTIVE: "
$part1 = "Eager to drive "
$part2 = "solutions at "
$part3 = "American Express"
$part4 = " on a full-time basis"

$pos1 = $pStart + $prefix.Length
$pos2 = $pos1 + $part1.Length
$pos3 = $pos2 + $part2.Length
$pos4 = $pos3 + $part3.Length
$pos5 = $pos4 + $part4.Length

$d.Range($pos1, $pos2).Font.Size = 12
$d.Range($pos2, $pos3).Font.Size = 12
$d.Range($pos3, $pos4).Font.Size = 12
$d.Range($pos4, $pos5).Font.Size = 12

# Re-insert the "_GoBack" bookmark between "Eager to drive " and
# "solutions at " (collapsed/zero-length), matching its new location.
$bmRange = $d.Range($pos2, $pos2)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
